$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 73 & 74: the two fixtures had been recorded in the wrong order; swap
# every column except A (the running index) between the two rows.
# ---------------------------------------------------------------------------
$ws.Cells.Item(73, 2).Value  = 7646749
$ws.Cells.Item(73, 5).Value  = "Brisbane Roar"
$ws.Cells.Item(73, 6).Value  = "Newcastle Jets"
$ws.Cells.Item(73, 8).Value  = 2
$ws.Cells.Item(73, 9).Value  = "H"
$ws.Cells.Item(73, 10).Value = 1.909
$ws.Cells.Item(73, 11).Value = 4
$ws.Cells.Item(73, 12).Value = 3.4
$ws.Cells.Item(73, 13).Value = 2.4
$ws.Cells.Item(73, 14).Value = 4
$ws.Cells.Item(73, 15).Value = 2.6
$ws.Cells.Item(73, 16).Value = 0
$ws.Cells.Item(73, 17).Value = 1.83
$ws.Cells.Item(73, 18).Value = 2.07
$ws.Cells.Item(73, 19).Value = 3.25
$ws.Cells.Item(73, 20).Value = 1.9
$ws.Cells.Item(73, 21).Value = 1.95
$ws.Cells.Item(73, 22).Value = 1.4
$ws.Cells.Item(73, 24).Value = -1
$ws.Cells.Item(73, 25).Value = 0.8300000000000001
$ws.Cells.Item(73, 26).Value = -1
$ws.Cells.Item(73, 27).Value = 0.8999999999999999

$ws.Cells.Item(74, 2).Value  = 7646750
$ws.Cells.Item(74, 5).Value  = "Perth Glory"
$ws.Cells.Item(74, 6).Value  = "Wellington Phoenix"
$ws.Cells.Item(74, 8).Value  = 4
$ws.Cells.Item(74, 9).Value  = "A"
$ws.Cells.Item(74, 10).Value = 2.45
$ws.Cells.Item(74, 11).Value = 3.75
$ws.Cells.Item(74, 12).Value = 2.55
$ws.Cells.Item(74, 13).Value = 3.1
$ws.Cells.Item(74, 14).Value = 3.8
$ws.Cells.Item(74, 15).Value = 2.05
$ws.Cells.Item(74, 16).Value = 0.25
$ws.Cells.Item(74, 17).Value = 2
$ws.Cells.Item(74, 18).Value = 1.85
$ws.Cells.Item(74, 19).Value = 3
$ws.Cells.Item(74, 20).Value = 1.925
$ws.Cells.Item(74, 21).Value = 1.925
$ws.Cells.Item(74, 22).Value = -1
$ws.Cells.Item(74, 24).Value = 1.05
$ws.Cells.Item(74, 25).Value = -1
$ws.Cells.Item(74, 26).Value = 0.8500000000000001
$ws.Cells.Item(74, 27).Value = 0.925

# ---------------------------------------------------------------------------
# Before touching rows 159/160 (whose B column holds the text fixtures
# "8109525" / "7127421" that later move down to rows 161/162), copy those
# two cells out verbatim (value-only paste preserves the text cell type
# without requiring any style/number-format change).
# ---------------------------------------------------------------------------
$ws.Cells.Item(159, 2).Copy()
$ws.Cells.Item(161, 2).PasteSpecial(-4163)
$ws.Cells.Item(160, 2).Copy()
$ws.Cells.Item(162, 2).PasteSpecial(-4163)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Rows 159 & 160 previously held placeholder fixtures (id only as text,
# missing score/odds columns). They now get the real recorded match data.
# ---------------------------------------------------------------------------
$ws.Cells.Item(159, 2).Value  = 7127419
$ws.Cells.Item(159, 3).Value  = "Australia ALeague"
$ws.Cells.Item(159, 4).Value  = 45409.17708333334
$ws.Cells.Item(159, 5).Value  = "Wellington Phoenix"
$ws.Cells.Item(159, 6).Value  = "Macarthur FC"
$ws.Cells.Item(159, 7).Value  = 3
$ws.Cells.Item(159, 8).Value  = 0
$ws.Cells.Item(159, 9).Value  = "H"
$ws.Cells.Item(159, 10).Value = 1.85
$ws.Cells.Item(159, 11).Value = 3.5
$ws.Cells.Item(159, 12).Value = 3.9
$ws.Cells.Item(159, 13).Value = 1.55
$ws.Cells.Item(159, 14).Value = 4.5
$ws.Cells.Item(159, 15).Value = 5.25
$ws.Cells.Item(159, 16).Value = -1
$ws.Cells.Item(159, 17).Value = 1.89
$ws.Cells.Item(159, 18).Value = 2.01
$ws.Cells.Item(159, 19).Value = 3.5
$ws.Cells.Item(159, 20).Value = 1.9
$ws.Cells.Item(159, 21).Value = 1.95
$ws.Cells.Item(159, 22).Value = 0.55
$ws.Cells.Item(159, 23).Value = -1
$ws.Cells.Item(159, 24).Value = -1
$ws.Cells.Item(159, 25).Value = 0.8899999999999999
$ws.Cells.Item(159, 26).Value = -1
$ws.Cells.Item(159, 27).Value = -1
$ws.Cells.Item(159, 28).Value = 0.95

$ws.Cells.Item(160, 2).Value  = 7127418
$ws.Cells.Item(160, 3).Value  = "Australia ALeague"
$ws.Cells.Item(160, 4).Value  = 45409.17708333334
$ws.Cells.Item(160, 5).Value  = "Newcastle Jets"
$ws.Cells.Item(160, 6).Value  = "Central Coast Mariners"
$ws.Cells.Item(160, 7).Value  = 1
$ws.Cells.Item(160, 8).Value  = 3
$ws.Cells.Item(160, 9).Value  = "A"
$ws.Cells.Item(160, 10).Value = 3.6
$ws.Cells.Item(160, 11).Value = 3.25
$ws.Cells.Item(160, 12).Value = 2
$ws.Cells.Item(160, 13).Value = 4.2
$ws.Cells.Item(160, 14).Value = 4
$ws.Cells.Item(160, 15).Value = 1.75
$ws.Cells.Item(160, 16).Value = 0.75
$ws.Cells.Item(160, 17).Value = 1.85
$ws.Cells.Item(160, 18).Value = 2
$ws.Cells.Item(160, 19).Value = 3
$ws.Cells.Item(160, 20).Value = 1.975
$ws.Cells.Item(160, 21).Value = 1.875
$ws.Cells.Item(160, 22).Value = -1
$ws.Cells.Item(160, 23).Value = -1
$ws.Cells.Item(160, 24).Value = 0.75
$ws.Cells.Item(160, 25).Value = -1
$ws.Cells.Item(160, 26).Value = 1
$ws.Cells.Item(160, 27).Value = 0.9750000000000001
$ws.Cells.Item(160, 28).Value = -1

# ---------------------------------------------------------------------------
# Rows 161 & 162 are brand-new rows holding the two still-unplayed fixtures
# that used to sit in rows 159/160 (odds refreshed for the upcoming games).
# Column A/D formatting is copied from row 160 so the new rows pick up the
# same bold/centered index style and date number format already in use.
# ---------------------------------------------------------------------------
$ws.Cells.Item(161, 1).Value = 159
$ws.Cells.Item(160, 1).Copy()
$ws.Cells.Item(161, 1).PasteSpecial(-4122)
$ws.Cells.Item(161, 3).Value = "Australia ALeague"
$ws.Cells.Item(161, 4).Value = 45410.08333333334
$ws.Cells.Item(160, 4).Copy()
$ws.Cells.Item(161, 4).PasteSpecial(-4122)
$ws.Cells.Item(161, 5).Value = "Sydney FC"
$ws.Cells.Item(161, 6).Value = "Perth Glory"
$ws.Cells.Item(161, 10).Value = 1.5
$ws.Cells.Item(161, 11).Value = 3.6
$ws.Cells.Item(161, 12).Value = 7
$ws.Cells.Item(161, 13).Value = 1.333
$ws.Cells.Item(161, 14).Value = 5.75
$ws.Cells.Item(161, 15).Value = 7.5
$ws.Cells.Item(161, 16).Value = -1.5
$ws.Cells.Item(161, 17).Value = 1.85
$ws.Cells.Item(161, 18).Value = 2.05
$ws.Cells.Item(161, 19).Value = 3.75
$ws.Cells.Item(161, 20).Value = 1.85
$ws.Cells.Item(161, 21).Value = 2
$ws.Cells.Item(161, 22).Value = 0
$ws.Cells.Item(161, 23).Value = 0
$ws.Cells.Item(161, 24).Value = 0

$ws.Cells.Item(162, 1).Value = 160
$ws.Cells.Item(160, 1).Copy()
$ws.Cells.Item(162, 1).PasteSpecial(-4122)
$ws.Cells.Item(162, 3).Value = "Australia ALeague"
$ws.Cells.Item(162, 4).Value = 45410.16666666666
$ws.Cells.Item(160, 4).Copy()
$ws.Cells.Item(162, 4).PasteSpecial(-4122)
$ws.Cells.Item(162, 5).Value = "Melbourne City"
$ws.Cells.Item(162, 6).Value = "Western United FC"
$ws.Cells.Item(162, 10).Value = 1.65
$ws.Cells.Item(162, 11).Value = 4
$ws.Cells.Item(162, 12).Value = 4.333
$ws.Cells.Item(162, 13).Value = 1.363
$ws.Cells.Item(162, 14).Value = 5.75
$ws.Cells.Item(162, 15).Value = 7
$ws.Cells.Item(162, 16).Value = -1.5
$ws.Cells.Item(162, 17).Value = 2.02
$ws.Cells.Item(162, 18).Value = 1.88
$ws.Cells.Item(162, 19).Value = 3.5
$ws.Cells.Item(162, 20).Value = 1.85
$ws.Cells.Item(162, 21).Value = 2
$ws.Cells.Item(162, 22).Value = 0
$ws.Cells.Item(162, 23).Value = 0
$ws.Cells.Item(162, 24).Value = 0

$excel.CutCopyMode = 0
